$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.57"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.45"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.730"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05833"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.417"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.469"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.313"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8014"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1467"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07625"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03217"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09248"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001670"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.249"

$ws.Range("E16").Value = "15MCDexMCBWorstin24h"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04760"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005975"

$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006266"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005425"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001065"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.699"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1264"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0009973"

$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04301"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007046"

$ws.Range("B42").Value = "BKEXToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1061"

$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"

$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003365"

$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009731"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005452"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7833"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.09943"

$ws.Range("E49").Value = "48BOLOBOLO"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002095"
